$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value, per the target diff.
$changes = @{
    24 = @{ 'J'=1.1; 'K'=7; 'L'=1.5; 'M'=2.5; 'N'=2.5; 'O'=1.5; 'P'=1.57;
            'Q'=2.25; 'R'=2.1; 'S'=1.67; 'AC'=81; 'AF'=8.5; 'AJ'=41 }
    37 = @{ 'G'=3; 'I'=2.63; 'J'=1.14; 'K'=5.5; 'L'=1.73; 'M'=2;
            'T'=6; 'U'=13; 'V'=13; 'W'=34; 'X'=34;
            'AE'=5.5; 'AF'=11; 'AG'=12; 'AH'=29; 'AI'=29 }
    44 = @{ 'G'=2.15; 'I'=3.1; 'N'=1.67; 'O'=2.15;
            'T'=10; 'U'=12; 'W'=21; 'Y'=21;
            'AF'=19; 'AG'=12; 'AH'=34; 'AI'=23; 'AJ'=29 }
    47 = @{ 'G'=1.8; 'H'=3.4; 'J'=1.07; 'K'=9; 'L'=1.33; 'M'=3.4; 'N'=2.08; 'O'=1.73;
            'AE'=12 }
    50 = @{ 'G'=1.73; 'H'=3.75; 'I'=4.75;
            'U'=7.5; 'W'=13; 'Z'=9.5;
            'AE'=12; 'AF'=23 }
    53 = @{ 'G'=2.2; 'H'=3.25; 'I'=3.25;
            'N'=2.5; 'O'=1.5;
            'T'=6; 'U'=9; 'V'=10; 'W'=21; 'X'=21; 'Z'=7;
            'AE'=7.5; 'AF'=15; 'AG'=13; 'AH'=41; 'AI'=34 }
    90 = @{ 'G'=2.45; 'H'=3.7; 'I'=2.63;
            'K'=17; 'L'=1.17; 'M'=5; 'N'=1.53; 'O'=2.4; 'P'=1.29; 'Q'=3.5; 'R'=1.5; 'S'=2.5;
            'T'=12; 'U'=15; 'V'=10; 'W'=26;
            'Z'=17; 'AA'=7.5; 'AB'=11; 'AD'=101;
            'AG'=10; 'AH'=26; 'AI'=19; 'AJ'=21 }
    91 = @{ 'G'=2.1; 'H'=3.5; 'I'=3.3;
            'N'=1.95; 'O'=1.9;
            'Z'=11; 'AC'=41; 'AE'=11 }
    150 = @{ 'G'=1.67; 'H'=3.7;
             'U'=7; 'X'=17;
             'AF'=23; 'AI'=41 }
    152 = @{ 'H'=3.5; 'J'=1.03; 'K'=15; 'N'=1.67; 'O'=2.15;
             'R'=1.57; 'S'=2.25; 'T'=10;
             'X'=15; 'Y'=21; 'AA'=7;
             'AC'=34; 'AD'=126; 'AE'=13; 'AF'=19; 'AJ'=26 }
    154 = @{ 'G'=2.7; 'H'=2.9; 'I'=2.7;
             'M'=2.37; 'N'=2.32; 'O'=1.47;
             'T'=6.9; 'U'=12.5; 'V'=10.25; 'W'=32; 'X'=26; 'Y'=40; 'Z'=6.6;
             'AB'=16.5;
             'AE'=6.6; 'AF'=12; 'AG'=10.5; 'AH'=32; 'AI'=28; 'AJ'=45 }
}

foreach ($rowNum in $changes.Keys) {
    $cols = $changes[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $val = $cols[$colLetter]
        $ws.Range("$colLetter$rowNum").Value = $val
    }
}
